# Pseudomonas putida LB-WT stationary phase reaction-sensitivity workbook:
# "results with fixed workflow" - the cutoff sweep now starts at cutoff-index 4
# instead of cutoff-index 0, so the per-cutoff Reaction_number series shifts up
# by 4 rows (the old rows 6..20 become the new rows 2..16) while the Cutoff
# column itself (A) keeps counting 0..14. The trailing 4 rows are dropped.

$wb = $excel.ActiveWorkbook

$sheetNames = @("NBR", "BAR")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Pull the old B/C values (rows 6..20) up into rows 2..16, leaving column A untouched.
    for ($r = 2; $r -le 16; $r++) {
        $srcRow = $r + 4
        $bVal = $ws.Cells.Item($srcRow, 2).Value2
        $cVal = $ws.Cells.Item($srcRow, 3).Value2
        $ws.Cells.Item($r, 2).Value = $bVal
        $ws.Cells.Item($r, 3).Value = $cVal
    }

    # Remove the now-duplicated trailing rows (17..20) and shrink the used range.
    $ws.Range("A17:C20").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
}
